# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the two
# outstanding files (2bb0d1bd-...md and 5a1b0d09-...md) have now been
# handed back and are in sync with en-US, for both the zh-cn and de-de
# locales:
#   - Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-locale sheets).
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     are populated on the per-locale sheets, each carrying its own
#     hyperlink (mirroring the existing Source File Name / Latest Handoff
#     File hyperlinks).
#   - The "Latest Handback DateTime" (G) column moves from the
#     "0001-01-01 00:00:00" placeholder to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrl1  = "https://github.com/OpenLocalizationTest/oltest/blob/79152dba4851d95f014590ff8b89d56dd9e14ef9/e2e/2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.md"
$mdUrl2  = "https://github.com/OpenLocalizationTest/oltest/blob/79152dba4851d95f014590ff8b89d56dd9e14ef9/e2e/5a1b0d09-fa51-4f90-b4f3-976d768860c4.md"

$zhXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd2772b5a6a9a30d1181a9abd67bfd36be38abf7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.56b2dcf44fd24a575f075277fbda96fee3a9ac9d.zh-cn.xlf"
$zhXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd2772b5a6a9a30d1181a9abd67bfd36be38abf7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5a1b0d09-fa51-4f90-b4f3-976d768860c4.f7341f2dcbdb90a5495014247d7315267f78aa28.zh-cn.xlf"

$deXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71a990df1679724d2192f5218eae80ed13a4115d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.56b2dcf44fd24a575f075277fbda96fee3a9ac9d.de-de.xlf"
$deXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71a990df1679724d2192f5218eae80ed13a4115d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5a1b0d09-fa51-4f90-b4f3-976d768860c4.f7341f2dcbdb90a5495014247d7315267f78aa28.de-de.xlf"

$mdName1 = "2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.md"
$mdName2 = "5a1b0d09-fa51-4f90-b4f3-976d768860c4.md"
$zhXlfName1 = "2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.56b2dcf44fd24a575f075277fbda96fee3a9ac9d.zh-cn.xlf"
$zhXlfName2 = "5a1b0d09-fa51-4f90-b4f3-976d768860c4.f7341f2dcbdb90a5495014247d7315267f78aa28.zh-cn.xlf"
$deXlfName1 = "2bb0d1bd-40f2-4d5d-92ee-29f4118fc308.56b2dcf44fd24a575f075277fbda96fee3a9ac9d.de-de.xlf"
$deXlfName2 = "5a1b0d09-fa51-4f90-b4f3-976d768860c4.f7341f2dcbdb90a5495014247d7315267f78aa28.de-de.xlf"

function Set-HandbackCell($ws, $addr, $text, $url, $display) {
    $ws.Range($addr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display) | Out-Null
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Overview sheet: the Status column (shared between zh-cn/de-de here)
# flips to the handback message for both tracked files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

Set-HandbackCell $wsZh "E2" $mdName1 $mdUrl1 $mdName1
Set-HandbackCell $wsZh "F2" $zhXlfName1 $zhXlfUrl1 $zhXlfName1
$wsZh.Range("G2").Value = "2016-03-08 02:41:42"

Set-HandbackCell $wsZh "E3" $mdName2 $mdUrl2 $mdName2
Set-HandbackCell $wsZh "F3" $zhXlfName2 $zhXlfUrl2 $zhXlfName2
$wsZh.Range("G3").Value = "2016-03-08 02:41:42"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

Set-HandbackCell $wsDe "E2" $mdName1 $mdUrl1 $mdName1
Set-HandbackCell $wsDe "F2" $deXlfName1 $deXlfUrl1 $deXlfName1
$wsDe.Range("G2").Value = "2016-03-08 02:41:57"

Set-HandbackCell $wsDe "E3" $mdName2 $mdUrl2 $mdName2
Set-HandbackCell $wsDe "F3" $deXlfName2 $deXlfUrl2 $deXlfName2
$wsDe.Range("G3").Value = "2016-03-08 02:41:57"
